$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Build the "Calibri 10, theme text color" look used by the new data row on a
# scratch cell well outside the used range, then fan it out via Copy /
# PasteSpecial (formats only) so every target cell shares the same style
# entry instead of minting a fresh one per cell.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.Font.Name = "Calibri"
$scratch.Font.ThemeColor = 1
$scratch.Copy()

$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)

[void]($excel.CutCopyMode = $false)
$scratch.Clear()

# ---------------------------------------------------------------------------
# New data row (row 2) — a Series-level entry for MCH341.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "MCH341-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33I | GRAP COUNT NUMER: NONE"

$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75

# ---------------------------------------------------------------------------
# View state: keep the header frozen, but move the live selection to E11.
# ---------------------------------------------------------------------------
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
[void]$ws.Range("E11").Select()
